$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old row 3 data (Specialist job)
$ws.Range("B2").Value = "Specialist – Operations, Fund Risk Management and Oversight (Open to all applicants)"
$ws.Range("D2").Value = "'02/04/2026"
$ws.Range("E2").Formula = '=HYPERLINK("https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/28371/?location=India&locationId=300000000440677&locationLevel=country&mode=location", "Apply")'

# Row 3 <- old row 4 data (Human Resources Associate)
$ws.Range("B3").Value = "Human Resources Associate"
$ws.Range("D3").Value = "'02/01/2026"
$ws.Range("E3").Formula = '=HYPERLINK("https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/31613/?location=India&locationId=300000000440677&locationLevel=country&mode=location", "Apply")'

# Row 4 keeps its Title/Location/Date (Human Resources Associate), only link changes
$ws.Range("E4").Formula = '=HYPERLINK("https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/31614/?location=India&locationId=300000000440677&locationLevel=country&mode=location", "Apply")'

# Row 5 <- old row 6 data (Finance Analyst)
$ws.Range("B5").Value = "Finance Analyst"
$ws.Range("D5").Value = "'01/26/2026"
$ws.Range("E5").Formula = '=HYPERLINK("https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/31468/?location=India&locationId=300000000440677&locationLevel=country&mode=location", "Apply")'

# Row 6 <- old row 7 data (Project Accounting & Financial Management Officer), date unchanged
$ws.Range("B6").Value = "Project Accounting & Financial Management Officer"
$ws.Range("E6").Formula = '=HYPERLINK("https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/31469/?location=India&locationId=300000000440677&locationLevel=country&mode=location", "Apply")'

# Remove the now-duplicate trailing row 7 and shrink used range
$ws.Rows("7").Delete()
